$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.240.37'
$ws.Range("E2").Value = '  +3.81%  '

$ws.Range("D3").Value = '2.621.19'
$ws.Range("E3").Value = '  +3.30%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.24'
$ws.Range("E5").Value = '  +2.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.07'
$ws.Range("E6").Value = '  +1.83%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  +0.92%  '

$ws.Range("D9").Value = '2.621.03'
$ws.Range("E9").Value = '  +3.30%  '

$ws.Range("E10").Value = '  +13.14%  '

$ws.Range("E11").Value = '  +0.16%  '

$ws.Range("E12").Value = '  +2.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.02'
$ws.Range("E13").Value = '  +0.14%  '

$ws.Range("D14").Value = '3.079.88'
$ws.Range("E14").Value = '  +3.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000186'
$ws.Range("E15").Value = '  +8.56%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.56'
$ws.Range("E16").Value = '  +1.64%  '

$ws.Range("D17").Value = '71.207.64'
$ws.Range("E17").Value = '  +3.92%  '

$ws.Range("D18").Value = '2.637.87'
$ws.Range("E18").Value = '  +3.89%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '384.01'
$ws.Range("E19").Value = '  +8.75%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.90'
$ws.Range("E20").Value = '  +5.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.48'
$ws.Range("E21").Value = '  +3.18%  '

$ws.Range("E22").Value = '  -1.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.26'
$ws.Range("E23").Value = '  +1.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.47'
$ws.Range("E24").Value = '  +5.62%  '

$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("E26").Value = '  +8.95%  '

$ws.Range("E27").Value = '  +5.59%  '

$ws.Range("D28").Value = '2.757.75'
$ws.Range("E28").Value = '  +4.50%  '

$ws.Range("E29").Value = '  -0.08%  '

$ws.Range("D30").Value = '0.0₃0963'
$ws.Range("E30").Value = '  +6.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '543.11'
$ws.Range("E31").Value = '  +5.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.05'
$ws.Range("E32").Value = '  +2.78%  '

$ws.Range("E33").Value = '  +4.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("E34").Value = '  +3.26%  '

$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.86'
$ws.Range("E36").Value = '  +1.72%  '

$ws.Range("E37").Value = '  -2.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.19'
$ws.Range("E38").Value = '  +4.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.89'
$ws.Range("E39").Value = '  +7.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.03'
$ws.Range("E40").Value = '  +1.74%  '

$ws.Range("E42").Value = '  +8.53%  '

$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.05'
$ws.Range("E44").Value = '  +3.86%  '

$ws.Range("E45").Value = '  +0.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.00'
$ws.Range("E46").Value = '  +2.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '154.28'
$ws.Range("E47").Value = '  +1.11%  '

$ws.Range("E48").Value = '  +1.80%  '

$ws.Range("E49").Value = '  +4.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.533'
$ws.Range("E50").Value = '  +2.22%  '

$ws.Range("D51").Value = '0.0₆0263'
$ws.Range("E51").Value = '  -0.17%  '
